$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Drop every existing hyperlink on the sheet up front -----------------
# NOTE: per-item Hyperlinks.Item(i).Delete() is a no-op in this runtime; only
# calling Delete() on the whole Hyperlinks collection actually removes them.
# We re-add the ones we still need (rows 2-11) further down.
$ws.Hyperlinks.Delete()

# --- 2) Overwrite rows 2-11 with the new listings ---------------------------
$ws.Range("A2").Value = "2025-12-16 06:30:30"
$ws.Range("B2").Value = "大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5455098"
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

$ws.Range("A3").Value = "2025-12-16 06:30:30"
$ws.Range("B3").Value = "【フルリモート】官公庁向けPythonアプリ開発PM募集|7名チーム統括"
$ws.Range("D3").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5454985"
$ws.Range("G3").Value = 295
$ws.Range("H3").Value = "🔥Python ◆開発 ◇アプリ"

$ws.Range("A4").Value = "2025-12-16 06:30:30"
$ws.Range("B4").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G4").Value = 243
$ws.Range("H4").Value = "🔥API ◆ツール"

$ws.Range("A5").Value = "2025-12-16 06:30:30"
$ws.Range("B5").Value = "【Java/対話システム/心理学実験】協同問題解決プラットフォームの改修開発"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5439921"
$ws.Range("G5").Value = 155
$ws.Range("H5").Value = "★Java ◆開発"

$ws.Range("A6").Value = "2025-12-16 06:30:30"
$ws.Range("B6").Value = "【急募】iPhone・Android対応の天気アプリ開発をお願いします!"
$ws.Range("D6").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5455038"
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = "◆開発 ◇アプリ"

$ws.Range("A7").Value = "2025-12-16 06:30:30"
$ws.Range("B7").Value = "ホームページ診断チェックツール"
$ws.Range("D7").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5455029"
$ws.Range("G7").Value = 73
$ws.Range("H7").Value = "◆ツール"

$ws.Range("A8").Value = "2025-12-16 06:30:30"
$ws.Range("B8").Value = "【急募】帳票デジタル化のフロントエンド開発者募集"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5454857"
$ws.Range("G8").Value = 68
$ws.Range("H8").Value = "◆開発"

$ws.Range("A9").Value = "2025-12-16 06:30:30"
$ws.Range("B9").Value = "【急募】Accessシステム改修・CSV読込・MySQLクラウド化・PDFデータ調整"
$ws.Range("D9").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5455015"
$ws.Range("G9").Value = 53
$ws.Range("H9").Value = "◇MySQL"

$ws.Range("A10").Value = "2025-12-16 06:30:30"
$ws.Range("B10").Value = "【急募】wixシステムでのメッセージ送信システム構築依頼"
$ws.Range("D10").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5455067"
$ws.Range("G10").Value = 33
$ws.Range("H10").ClearContents()

$ws.Range("A11").Value = "2025-12-16 06:30:30"
$ws.Range("B11").Value = "【SESエンジニア募集】多様なプロジェクトに参画可能!"
$ws.Range("D11").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5437544"
$ws.Range("G11").Value = 25
$ws.Range("H11").ClearContents()

# --- 3) Drop the now-stale rows 12-21 entirely ------------------------------
$ws.Rows("12:21").Delete()

# --- 4) Re-create the F2:F11 hyperlinks against the refreshed URLs ---------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5455098")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5454985")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5439921")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5455038")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5455029")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5454857")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5455015")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5455067")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5437544")

# --- 5) Column width tweaks (col B: 52 -> 47, col H: 19 -> 18) -------------
# ColumnWidth as read/written through this COM layer is offset by ~0.8333
# (5/6) character-widths from the raw OOXML <col width> value, so we dial
# the COM-side figure back by that much to land on the exact OOXML target.
$ws.Columns.Item(2).ColumnWidth = 46.16666666666666
$ws.Columns.Item(8).ColumnWidth = 17.16666666666667
